$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A83:E83").NumberFormat = "@"

$ws.Range("A83").Value = "2025-12-08"
$ws.Range("B83").Value = "Pick 3"
$ws.Range("C83").Value = "251208"
$ws.Range("D83").Value = "4-1-4"
$ws.Range("E83").Value = "2025-12-08T21:42:07.627+04:00"
